# Update "想去人数" (F column) values in the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 599
$ws1.Range("F3").Value = 202
$ws1.Range("F4").Value = 463
$ws1.Range("F5").Value = 490
$ws1.Range("F7").Value = 2538
$ws1.Range("F8").Value = 432
$ws1.Range("F9").Value = 6837
$ws1.Range("F10").Value = 182
$ws1.Range("F11").Value = 433
$ws1.Range("F12").Value = 4

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 599
$ws4.Range("F3").Value = 202
$ws4.Range("F4").Value = 463
$ws4.Range("F5").Value = 490
$ws4.Range("F9").Value = 2538
$ws4.Range("F10").Value = 432
$ws4.Range("F11").Value = 6837
$ws4.Range("F12").Value = 182
$ws4.Range("F13").Value = 433
$ws4.Range("F14").Value = 4
